$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.09039999999999999
$ws.Range("B3").Value = 0.0804
$ws.Range("B4").Value = 0.0216
$ws.Range("B5").Value = 0.15
$ws.Range("B6").Value = 0.1436
$ws.Range("B7").Value = 0.11
$ws.Range("B8").Value = 0.0256
$ws.Range("B9").Value = 0.0532
$ws.Range("B10").Value = 0.08400000000000001
$ws.Range("B11").Value = 0.04960000000000001
$ws.Range("B12").Value = 0.052
$ws.Range("B13").Value = 0.326
$ws.Range("B14").Value = 0.0588
$ws.Range("B15").Value = 0.0212
$ws.Range("B16").Value = 0.1976
$ws.Range("B17").Value = 0.0484
$ws.Range("B18").Value = 0.1124
$ws.Range("B19").Value = 0.152
$ws.Range("B20").Value = 0.1892
$ws.Range("B21").Value = 0.0132
$ws.Range("B22").Value = 0.2128
$ws.Range("B23").Value = 0.0672
$ws.Range("B24").Value = 0.06319999999999999
$ws.Range("B25").Value = 0.1264
$ws.Range("B26").Value = 0.0328
$ws.Range("B27").Value = 0.2528
$ws.Range("B28").Value = 0.1692
$ws.Range("B29").Value = 0.022
$ws.Range("B30").Value = 0.0212
$ws.Range("B31").Value = 0.01
$ws.Range("B32").Value = 0.008399999999999999
$ws.Range("B33").Value = 0.1516
$ws.Range("B34").Value = 0.1404
$ws.Range("B35").Value = 0.0412
$ws.Range("B36").Value = 0.198
$ws.Range("B37").Value = 0.1652
$ws.Range("B38").Value = 0.1932
$ws.Range("B39").Value = 0.0252
$ws.Range("B40").Value = 0.106
$ws.Range("B41").Value = 0.1276
$ws.Range("B42").Value = 0.0588
$ws.Range("B43").Value = 0.0136
$ws.Range("B44").Value = 0.2084
$ws.Range("B45").Value = 0.0224
$ws.Range("B46").Value = 0.0152
$ws.Range("B47").Value = 0.0192
$ws.Range("B48").Value = 0.06279999999999999
$ws.Range("B49").Value = 0.0832
$ws.Range("B50").Value = 0.0336
$ws.Range("B51").Value = 0.008399999999999999
$ws.Range("B52").Value = 0.0956
$ws.Range("B53").Value = 0.256
$ws.Range("B54").Value = 0.0344
$ws.Range("B55").Value = 0.02
$ws.Range("B56").Value = 0.0488
$ws.Range("B57").Value = 0.1356
$ws.Range("B58").Value = 0.04
$ws.Range("B59").Value = 0.0372
$ws.Range("B60").Value = 0.0304
$ws.Range("B61").Value = 0.0112
$ws.Range("B62").Value = 0.0492
$ws.Range("B63").Value = 0.07919999999999999
$ws.Range("B64").Value = 0.0316
$ws.Range("B65").Value = 0.0548
$ws.Range("B66").Value = 0.0572
$ws.Range("B67").Value = 0.1068
$ws.Range("B68").Value = 0.7396
$ws.Range("B69").Value = 0.0196
$ws.Range("B70").Value = 0.08800000000000001
$ws.Range("B71").Value = 0.0144
$ws.Range("B72").Value = 0.0616
$ws.Range("B73").Value = 0.06999999999999999
$ws.Range("B74").Value = 0.0832
$ws.Range("B75").Value = 0.0092
$ws.Range("B76").Value = 0.0068
$ws.Range("B77").Value = 0.04960000000000001
$ws.Range("B78").Value = 0.0316
$ws.Range("B79").Value = 0.9608
$ws.Range("B80").Value = 0.0412
$ws.Range("B81").Value = 0.1332
$ws.Range("B82").Value = 0.1536
$ws.Range("B83").Value = 0.0716
$ws.Range("B84").Value = 0.0512
$ws.Range("B85").Value = 0.0308
$ws.Range("B86").Value = 0.0668
$ws.Range("B87").Value = 0.194
$ws.Range("B88").Value = 0.112
$ws.Range("B89").Value = 0.0196
$ws.Range("B90").Value = 0.022
$ws.Range("B91").Value = 0.08639999999999999
$ws.Range("B92").Value = 0.2164
$ws.Range("B93").Value = 0.0052
$ws.Range("B94").Value = 0.0616
$ws.Range("B95").Value = 0.136
$ws.Range("B96").Value = 0.2552
$ws.Range("B97").Value = 0.03759999999999999
$ws.Range("B98").Value = 0.0092
$ws.Range("B99").Value = 0.0148
$ws.Range("B100").Value = 0.1188
$ws.Range("B101").Value = 0.0612
$ws.Range("B102").Value = 0.0144
$ws.Range("B103").Value = 0.008
$ws.Range("B104").Value = 0.1532
$ws.Range("B105").Value = 0.2632
$ws.Range("B106").Value = 0.0164
$ws.Range("B107").Value = 0.0428
$ws.Range("B108").Value = 0.1444
$ws.Range("B109").Value = 0.1652
$ws.Range("B110").Value = 0.0716
$ws.Range("B111").Value = 0.1032
$ws.Range("B112").Value = 0.2948
$ws.Range("B113").Value = 0.1232
$ws.Range("B114").Value = 0.008
$ws.Range("B115").Value = 0.004
$ws.Range("B116").Value = 0.0172
$ws.Range("B117").Value = 0.0448
$ws.Range("B118").Value = 0.1528
$ws.Range("B119").Value = 0.052
$ws.Range("B120").Value = 0.1012
$ws.Range("B121").Value = 0.0372
$ws.Range("B122").Value = 0.1324
$ws.Range("B123").Value = 0.2296
$ws.Range("B124").Value = 0.3144
$ws.Range("B125").Value = 0.048
$ws.Range("B126").Value = 0.1396
$ws.Range("B127").Value = 0.1024
$ws.Range("B128").Value = 0.0104
$ws.Range("B129").Value = 0.0148
$ws.Range("B130").Value = 0.0872
$ws.Range("B131").Value = 0.006399999999999999
$ws.Range("B132").Value = 0.04
$ws.Range("B133").Value = 0.0616
$ws.Range("B134").Value = 0.0772
$ws.Range("B135").Value = 0.0696
$ws.Range("B136").Value = 0.0224
$ws.Range("B137").Value = 0.098
$ws.Range("B138").Value = 0.0204
$ws.Range("B139").Value = 0.0776
$ws.Range("B140").Value = 0.1104
$ws.Range("B141").Value = 0.1896
$ws.Range("B142").Value = 0.0504
$ws.Range("B143").Value = 0.0108
$ws.Range("B144").Value = 0.1328
$ws.Range("B145").Value = 0.0384
$ws.Range("B146").Value = 0.0868
$ws.Range("B147").Value = 0.0092
$ws.Range("B148").Value = 0.0052
$ws.Range("B149").Value = 0.0256
$ws.Range("B150").Value = 0.0176
$ws.Range("B151").Value = 0.07479999999999999
$ws.Range("B152").Value = 0.0664
$ws.Range("B153").Value = 0.5232
$ws.Range("B154").Value = 0.182
$ws.Range("B155").Value = 0.0232
$ws.Range("B156").Value = 0.022
$ws.Range("B157").Value = 0.0848
$ws.Range("B158").Value = 0.124
$ws.Range("B159").Value = 0.1812
$ws.Range("B160").Value = 0.08800000000000001
$ws.Range("B161").Value = 0.0152
$ws.Range("B162").Value = 0.0348
$ws.Range("B163").Value = 0.0268
$ws.Range("B164").Value = 0.1228
$ws.Range("B165").Value = 0.0692
$ws.Range("B166").Value = 0.0544
$ws.Range("B167").Value = 0.02
$ws.Range("B168").Value = 0.2264
$ws.Range("B169").Value = 0.0196
$ws.Range("B170").Value = 0.1188
$ws.Range("B171").Value = 0.0216
$ws.Range("B172").Value = 0.0048
$ws.Range("B173").Value = 0.08599999999999999
$ws.Range("B174").Value = 0.0328
$ws.Range("B175").Value = 0.1732
$ws.Range("B176").Value = 0.0544
$ws.Range("B177").Value = 0.028
$ws.Range("B178").Value = 0.0152
$ws.Range("B179").Value = 0.0568
$ws.Range("B180").Value = 0.0776
$ws.Range("B181").Value = 0.1416
$ws.Range("B182").Value = 0.0436
$ws.Range("B183").Value = 0.3912
